# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet and
#    populate it with the fund-holding detail rows.
# 2. Insert a new top data row in "总计" summarizing the new quarter and
#    renumber the existing index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet
# ---------------------------------------------------------------------
$totalSheetAnchor = $wb.Worksheets.Item("总计")

$ws = $wb.Worksheets.Add($totalSheetAnchor)
$ws.Name = "2022-Q1"

# NOTE: the worksheet object passed as the "Before" anchor above becomes
# stale once the new sheet is inserted - always re-resolve sheets by
# name after calling Worksheets.Add().
$totalSheet = $wb.Worksheets.Item("总计")
$q4Sheet = $wb.Worksheets.Item("2021-Q4")

# Header row text
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Data rows - B..G hold numeric-looking text, so force a text number
# format before writing, then clear the format again afterwards so the
# cells end up styled exactly like brand-new, never-formatted cells
# (matching the source workbook's minimal styling).
$textRange = $ws.Range("B2:G8")
$textRange.NumberFormat = "@"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "910021"
$ws.Range("C2").Value = "东方红启华三年持有期混合型证券投资基金A"
$ws.Range("D2").Value = "5.93"
$ws.Range("E2").Value = "86.11"
$ws.Range("F2").Value = "3.16"
$ws.Range("G2").Value = "0.1874"
$ws.Range("H2").Value = 7

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "690001"
$ws.Range("C3").Value = "民生加银品牌蓝筹混合"
$ws.Range("D3").Value = "1.21"
$ws.Range("E3").Value = "93.01"
$ws.Range("F3").Value = "4.04"
$ws.Range("G3").Value = "0.0489"
$ws.Range("H3").Value = 8

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "011313"
$ws.Range("C4").Value = "东方红启华三年持有期混合型证券投资基金B"
$ws.Range("D4").Value = "0.97"
$ws.Range("E4").Value = "86.11"
$ws.Range("F4").Value = "3.16"
$ws.Range("G4").Value = "0.0307"
$ws.Range("H4").Value = 7

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "970048"
$ws.Range("C5").Value = "东海证券海睿致远灵活配置混合型集合资产管理计划"
$ws.Range("D5").Value = "0.54"
$ws.Range("E5").Value = "79.71"
$ws.Range("F5").Value = "3.65"
$ws.Range("G5").Value = "0.0197"
$ws.Range("H5").Value = 7

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "620004"
$ws.Range("C6").Value = "金元顺安价值增长混合"
$ws.Range("D6").Value = "0.18"
$ws.Range("E6").Value = "87.28"
$ws.Range("F6").Value = "3.13"
$ws.Range("G6").Value = "0.0056"
$ws.Range("H6").Value = 7

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "620002"
$ws.Range("C7").Value = "金元顺安成长动力混合"
$ws.Range("D7").Value = "0.16"
$ws.Range("E7").Value = "62.78"
$ws.Range("F7").Value = "3.09"
$ws.Range("G7").Value = "0.0049"
$ws.Range("H7").Value = 5

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "005021"
$ws.Range("C8").Value = "渤海汇金量化汇盈灵活配置混合"
$ws.Range("D8").Value = "0.02"
$ws.Range("E8").Value = "92.66"
$ws.Range("F8").Value = "1.12"
$ws.Range("G8").Value = "0.0002"
$ws.Range("H8").Value = 10

$textRange.ClearFormats()

# Re-apply the shared "header / index" look (bold, thin border, centered)
# by copying the formatting already used by the same cells on the
# "2021-Q4" sheet, which is laid out identically.
$q4Sheet.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$q4Sheet.Range("A2:A4").Copy()
$ws.Range("A2:A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Prepend the new quarter to the "总计" summary sheet
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 0.3

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
